$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 809, pushing existing rows 809:850 down to 810:851
$ws.Rows("809:809").Insert()

# Populate the newly inserted row with the new record.
# Force Text format on A so the date-like string isn't auto-converted
# to a date serial number (matches the existing text date cells), then
# reset the cell style back to Normal so no stray style index is left
# behind (the source rows carry no explicit style).
$ws.Range("A809").NumberFormat = "@"
$ws.Range("A809").Value = "2026/02/17"
$ws.Range("A809").Style = "Normal"
$ws.Range("B809").Value = "火"
$ws.Range("C809").Value = 10
$ws.Range("D809").Value = 58
